$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Each value is prefixed with a
# leading apostrophe so Excel stores it as literal text (matching the
# original workbook, where these cells hold text such as "20.526.92" or
# "  +1.47%  " rather than numbers/dates/percentages).
$updates = @{
    "D2" = '''20.526.92'
    "E2" = '''  +1.47%  '
    "D3" = '''1.472.57'
    "E3" = '''  +1.98%  '
    "D4" = '''1.008'
    "E4" = '''  +0.19%  '
    "D5" = '''0.9569'
    "E5" = '''  +4.87%  '
    "D6" = '''277.59'
    "E6" = '''  +0.03%  '
    "D7" = '''0.3614'
    "E7" = '''  -1.40%  '
    "D8" = '''0.3064'
    "E8" = '''  -2.32%  '
    "D9" = '''39.58'
    "E9" = '''  +1.84%  '
    "D10" = '''1.063'
    "E10" = '''  +4.06%  '
    "D11" = '''0.06635'
    "E11" = '''  +1.51%  '
    "D12" = '''1.002'
    "E12" = '''  +0.19%  '
    "D13" = '''5.524'
    "E13" = '''  +2.29%  '
    "D14" = '''18.08'
    "E14" = '''  +2.92%  '
    "D15" = '''6.183'
    "E15" = '''  +1.73%  '
    "D16" = '''0.9573'
    "E16" = '''  +2.47%  '
    "D17" = '''0.00001025'
    "E17" = '''  +0.85%  '
    "D18" = '''1.475.52'
    "E18" = '''  +2.40%  '
    "D19" = '''0.05916'
    "E19" = '''  +5.02%  '
    "D20" = '''69.02'
    "E20" = '''  +1.79%  '
    "D21" = '''5.494'
    "E21" = '''  +1.61%  '
    "D22" = '''14.49'
    "E22" = '''  -0.02%  '
    "D23" = '''11.17'
    "E23" = '''  +3.25%  '
    "D24" = '''2.254'
    "E24" = '''  +0.11%  '
    "D25" = '''20.548.74'
    "E25" = '''  +1.68%  '
    "D26" = '''143.58'
    "E26" = '''  +4.98%  '
    "D27" = '''2.120'
    "E27" = '''  -3.16%  '
    "D28" = '''17.16'
    "E28" = '''  +1.20%  '
    "D29" = '''1.637.43'
    "E29" = '''  +2.70%  '
    "D30" = '''113.63'
    "E30" = '''  +2.89%  '
    "D31" = '''3.911'
    "E31" = '''  +3.25%  '
    "D32" = '''4.964'
    "E32" = '''  +2.54%  '
    "B33" = '''Stellar'
    "C33" = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "D33" = '''0.07994'
    "E33" = '''  +3.83%  '
    "B34" = '''ImmutableX'
    "C34" = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "D34" = '''0.8073'
    "E34" = '''  +0.11%  '
    "D35" = '''1.515'
    "E35" = '''  +4.43%  '
    "D36" = '''1.216'
    "E36" = '''  +6.75%  '
    "D37" = '''0.05758'
    "E37" = '''  -4.12%  '
    "D38" = '''4.723'
    "E38" = '''  +0.37%  '
    "D39" = '''0.02054'
    "E39" = '''  +2.82%  '
    "B40" = '''Frax'
    "C40" = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    "D40" = '''0.9579'
    "E40" = '''  +2.76%  '
    "B41" = '''Aptos'
    "C41" = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    "D41" = '''10.36'
    "E41" = '''  +1.48%  '
    "D42" = '''0.1876'
    "E42" = '''  +2.15%  '
    "D43" = '''7.429'
    "E43" = '''  +2.76%  '
    "D44" = '''0.5268'
    "E44" = '''  +0.45%  '
    "D45" = '''12.26'
    "E45" = '''  +2.16%  '
    "D46" = '''3.525'
    "E46" = '''  +0.01%  '
    "D47" = '''118.00'
    "E47" = '''  -0.54%  '
    "D48" = '''0.5198'
    "E48" = '''  +0.98%  '
    "D49" = '''1.813'
    "E49" = '''  +2.53%  '
    "D50" = '''0.06467'
    "E50" = '''  +2.19%  '
    "D51" = '''0.9863'
    "E51" = '''  -0.86%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
